$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update raw input values in column C (counts), D has formulas C#/C5 that recalc automatically
$ws.Range("C3").Value = 238
$ws.Range("C4").Value = 467
$ws.Range("C5").Value = 542.9
$ws.Range("C7").Value = 1129
$ws.Range("C8").Value = 2572
$ws.Range("C9").Value = 2045

# C3 and C4 lose their custom number-format style (back to default/General)
$ws.Range("C3").Style = "Normal"
$ws.Range("C4").Style = "Normal"

# Update the selection to C6 (matches the author's last-edited cell)
[void]$ws.Range("C6").Select()

[void]$excel.CalculateFull()
